$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.274.11"
$ws.Range("E2").Value = "  +1.77%  "
$ws.Range("D3").Value = "3.736.25"
$ws.Range("E3").Value = "  +0.85%  "
$ws.Range("E4").Value = "  -0.64%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "613.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.97%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "187.78"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.84%  "
$ws.Range("D7").Value = "3.732.83"
$ws.Range("E7").Value = "  +0.66%  "
$ws.Range("E8").Value = "  +1.10%  "
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.723"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.82%  "
$ws.Range("E11").Value = "  -2.33%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "57.16"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +6.85%  "
$ws.Range("E13").Value = "  -1.78%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.70"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.31%  "
$ws.Range("D15").Value = "4.331.57"
$ws.Range("E15").Value = "  +0.61%  "
$ws.Range("D16").Value = "3.731.55"
$ws.Range("E16").Value = "  +0.11%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.11"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.16%  "
$ws.Range("E18").Value = "  -0.17%  "
$ws.Range("E19").Value = "  -0.36%  "
$ws.Range("E20").Value = "  -0.14%  "
$ws.Range("D21").Value = "69.059.45"
$ws.Range("E21").Value = "  +1.19%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "414.72"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.93%  "
$ws.Range("E23").Value = "  +1.32%  "
$ws.Range("E24").Value = "  +0.08%  "
$ws.Range("E25").Value = "  -0.23%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.94"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.29%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.99"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.30%  "
$ws.Range("E28").Value = "  +2.30%  "
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.70"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.75%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "33.46"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.08%  "
$ws.Range("E32").Value = "  -11.16%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "12.79"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("E34").Value = "  +3.44%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "44.90"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.27%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "624.38"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.39%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "65.69"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.45%  "
$ws.Range("D38").Value = "0.0₃0850"
$ws.Range("E38").Value = "  -9.40%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.414"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.82%  "
$ws.Range("E40").Value = "  +0.27%  "
$ws.Range("E41").Value = "  -0.98%  "
$ws.Range("E42").Value = "  +3.08%  "
$ws.Range("E43").Value = "  -0.68%  "
$ws.Range("E44").Value = "  +0.87%  "
$ws.Range("E45").Value = "  +1.39%  "
$ws.Range("E46").Value = "  +4.52%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.23"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.29%  "
$ws.Range("D48").Value = "2.837.10"
$ws.Range("E48").Value = "  +3.14%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.76"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.49%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.72"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -16.67%  "
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "143.82"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.01%  "
